# Edit script for dig_information.xlsx
# Commit message: replaced mathjax with svg files in model landing pages
#
# Summary of changes:
# 1. Active sheet changes from "params_deep" to "equations" (workbook activeTab 3 -> 5)
# 2. "params_deep" sheet view: selection moves from A34 to E6 (no longer the active tab)
# 3. "equations" sheet: duplicate row 2/3 (both referred to the same two strings) is
#    fixed by deleting the stray row 3; rows below shift up. Column A becomes a running
#    count formula (=A{n-1}+1) instead of a hard-coded number, and the selection on
#    that sheet moves to B45 (last data cell), becoming the active tab.
# 4. The LaTeX equation strings in column C have their escaped double backslashes
#    ("\\") collapsed to single backslashes ("\"), and a few strings get their
#    inline "\ \ ... \ " spacing hacks replaced with "\quad \text{...} \quad ...".

$wb = $excel.ActiveWorkbook

$wsParamsDeep = $wb.Worksheets.Item("params_deep")
$wsEquations = $wb.Worksheets.Item("equations")

# --- Fix the "equations" sheet: remove the duplicated row, renumber column A ---

# Row 3 duplicates row 2 (same two shared strings). Deleting it shifts every
# subsequent row up by one, which both removes the duplicate pair and makes
# each following row line up with the correct (next) equation text.
$wsEquations.Rows.Item(3).Delete()

# Column A used to hold static sequence numbers (0,1,2,...,44). After the
# deletion it should become a running-count formula referencing the row above.
for ($r = 3; $r -le 45; $r++) {
    $prevRow = $r - 1
    $wsEquations.Range("A$r").Formula = "=A$prevRow+1"
}

# The equation text (column C) had its LaTeX strings double-escaped
# (e.g. "\\left" instead of "\left"). Replace each with the corrected text.
$equationTexts = @(
    'q_{x,t}=A_{x,t}\left({z^e}_{t-1}\right)^{\psi_{x}}\left(k_{x,t-1}\right)^{\alpha_{x}}\left(L_{x,t}\right)^{1-\alpha_{x}}',
    'q_{n,t}=A_{n,t}\left({z^e}_{t-1}\right) ^{\psi_{n}}\left( k_{n,t-1}\right)^{\alpha_{n}}\left( L_{n,t}\right)^{1-\alpha_{n}}',
    'A_{x,t}=a_{x}\left( \frac{{q^{I}}_{x,t-1}}{\bar{q^{I}}_{x}}\right) ^{\sigma_{x}}\left( {k^{I}}_{x,t-1}\right) ^{\xi_{x}}',
    'A_{n,t}=a_{n}\left( \frac{{q^{I}}_{n,t-1}}{\bar{q^{I}}_{n}}\right)^{\sigma_{n}}\left( {k^{I}}_{n,t-1}\right) ^{\xi_{n}}',
    'P_{k,t}=P_{mm,t}+a_{k}P_{n,t}',
    'P_{z,t}=P_{mm,t}+a_{z}P_{n,t}',
    'P_{n,t}(1-\alpha_{n})\frac{q_{n,t}}{L_{n,t}}=w_{t}',
    'P_{x,t}(1-\alpha_{x})\frac{q_{x,t}}{L_{x,t}}=w_{t}',
    'P_{n,t}\alpha_{n}\frac{q_{n,t}}{k_{n,t-1}}=r_{n,t}',
    'P_{x,t}\alpha_{x}\frac{q_{x,t}}{k_{x,t-1}}=r_{x,t}',
    '{c^{i}}_{t}=\left[ \rho_{x}^{\frac{1}{\epsilon }}\left( {c^{i}}_{x,t}\right) ^{\frac{\epsilon-1}{\epsilon}}+\rho_{m}^{\frac{1}{\epsilon }}\left({c^{i}}_{m,t}\right) ^{\frac{\epsilon-1}{\epsilon}}+(\rho_{n})^{\frac{1}{\epsilon }}\left( {c^{i}}_{n,t}\right) ^{\frac{\epsilon -1}{\epsilon }}\right]^{^{\frac{\epsilon }{\epsilon -1}}}\quad \text{for} \quad i=\mathfrak{s},\mathfrak{h}',
    'P_{t}=\left[ \rho_{x}P_{x,t}^{1-\epsilon }+\rho_{m}P_{m,t}^{1-\epsilon }+\rho_{n}P_{n,t}^{1-\epsilon }\right] ^{^{\frac{1}{1-\epsilon }}}',
    '{c^{i}}_{j,t}=\rho _{j}\left( \frac{P_{j,t}}{P_{t}}\right) ^{-\epsilon} {c^{i}}_{t} \quad \text{for} \quad j=x,m,n \quad \text{and} \quad i=\mathfrak{s},\mathfrak{h}',
    'P_{t}{b^{\mathfrak{s}}}_{t}-{b^{\mathfrak{s}\ast}}_{t} = r_{x,t} {k^{\mathfrak{s}}}_{x,t-1} + r_{n,t-1} {k^{\mathfrak{s}}}_{n,t-1} + w_{t}L_{t}^{\mathfrak{s}}+ \frac{\mathcal{R}_{t}}{1+a} +\frac{\mathcal{T}_{t}}{1+a}-\frac{1+{r^{\ast}}_{t-1}}{1+g} {b^{\mathfrak{s}\ast}}_{t-1} +\frac{1+r_{t-1}}{1+g} P_{t} {b^{\mathfrak{s}}}_{t-1} - P_{k,t}\left( {i^{\mathfrak{s}}}_{x,t} +{i^{\mathfrak{s}}}_{n,t} + {AC^{\mathfrak{s}}}_{x,t} + {AC^{\mathfrak{s}}}_{n,t}\right) - P_{t}{c^{\mathfrak{s}}_{t}}(1+h_{t})-\mu {z^{e}}_{t-1}-\mathcal{P^{\mathfrak{s}}}_{t} - {\Phi^{\mathfrak{s}}}_{t}',
    '(1+g) {k^{\mathfrak{s}}}_{x,t} = {i^{\mathfrak{s}}}_{x,t}+(1-\delta_{x}) {k^{\mathfrak{s}}}_{x,t-1}',
    '(1+g) {k^{\mathfrak{s}}}_{n,t} = {i^{\mathfrak{s}}}_{n,t}+(1-\delta_{n}) {k^{\mathfrak{s}}}_{n,t-1}',
    '{AC^{\mathfrak{s}}}_{j,t} \equiv \frac{v}{2}\left( \frac{{i^{\mathfrak{s}}}_{j,t}}{{k^{\mathfrak{s}}}_{j,t-1}}-\delta_{j} -g\right) ^{2} {k^{\mathfrak{s}}}_{j,t-1} \quad \text{for} \quad j=x,n',
    '\mathcal{P^{\mathfrak{s}}}_{t} \equiv \frac{\eta}{2}({b^{\mathfrak{s}\ast}}_{t} - \bar{b^{\mathfrak{s}\ast}})^{2}',
    '{c^{\mathfrak{s}}}_{t}={c^{\mathfrak{s}}}_{t+1}\left( \beta \frac{1+r_{t}}{1+g}\frac{1+h_{t}}{1+h_{t+1}}\right) ^{-\tau }',
    '(1+r_{t})\frac{P_{t+1}}{P_{t}}=\frac{1+{r^{\ast}}_{t}}{\left[ 1-\eta ({b^{\mathfrak{s}\ast}}_{t}-\bar{b^{\mathfrak{s}\ast}})\right] }',
    '\frac{r_{x,t+1}}{P_{k,t+1}}+1-\delta_{x} +v{\Upsilon^{\mathfrak{s}}}_{x,t+1} \left( \frac{{i^{\mathfrak{s}}}_{x,t+1}}{{k^{\mathfrak{s}}}_{x,t}}+1-\delta_{x}\right) -\frac{v}{2}\left( {\Upsilon^{\mathfrak{s}}}_{x,t+1}\right)^{2} =(1+r_{t})\frac{P_{t+1}}{P_{t}}\frac{P_{k,t}}{P_{k,t+1}}\left(1+v{\Upsilon^{\mathfrak{s}}}_{x,t}\right)',
    '\frac{r_{n,t+1}}{P_{k,t+1}}+1-\delta_{n} +v{\Upsilon^{\mathfrak{s}}}_{n,t+1} \left( \frac{{i^{\mathfrak{s}}}_{n,t+1}}{{k^{\mathfrak{s}}}_{n,t}}+1-\delta_{n}\right) -\frac{v}{2}\left( {\Upsilon^{\mathfrak{s}}}_{n,t+1}\right)^{2} = (1+r_{t})\frac{P_{t+1}}{P_{t}}\frac{P_{k,t}}{P_{k,t+1}}\left(1+v{\Upsilon^{\mathfrak{s}}}_{n,t}\right)',
    '{\Upsilon^{\mathfrak{s}}}_{j,t}=\left( \frac{{i^{\mathfrak{s}}}_{j,t}}{{k^{\mathfrak{s}}}_{j,t-1}}-\delta_{j} -g\right) \quad \text{for}\quad j=x,n',
    '\eta ({b^{\mathfrak{s}\ast}}_{t}-\bar{b}^{\mathfrak{s}\ast})=1-\frac{1+{r^{\ast}}_{t}}{(1+r_{t})\frac{P_{t+1}}{P_{t}}}',
    '{r^{\ast}}_{t}=r_{dc,t}+\mathfrak{u}',
    '(1+h_{t})P_{t} {c^{\mathfrak{h}}}_{t}=w_{t}L^{\mathfrak{h}}+\frac{a}{1+a}(\mathcal{R}_{t}+\mathcal{T}_{t})',
    '(1+g)z_{t}=(1-\delta_{z} )z_{t-1}+i_{z,t}',
    'z_{t}^{e}=\bar{s}\bar{z}+s(z_{t}-\bar{z}) \quad \text{with} \quad \bar{s} \in \lbrack 0,1] \quad \text{and} \quad s \in \lbrack 0,1]',
    '(1+g) {z^{e}}_{t}=(1-\delta_{z} ){z^{e}}_{t-1}+s(i_{z,t}-\bar{\imath}_{z})+\bar{s}\bar{\imath}_{z}',
    'P_{t}\Delta b_{t}+\Delta d_{c,t}+\Delta d_{t} = \frac{r_{t-1}-g}{1+g}P_{t}b_{t-1}+\frac{r_{d,t-1}-g}{1+g}d_{t-1}+\frac{r_{dc,t-1}-g}{1+g}d_{c,t-1} + P_{z,t}\mathbb{I}_{z,t}+\mathcal{T}_{t}-h_{t}P_{t}c_{t}-\mathcal{G}_{t}-\mathcal{N}_{t}-\mu {z^{e}}_{t-1}',
    'r_{dc,t}=r^{f}+\upsilon _{g}e^{\eta _{g}\left( \frac{d_{t}+d_{c,t}}{y_{t}}-\frac{\bar{d}+\bar{d}_{c}}{\bar{y}}\right) }',
    '\mathbb{I}_{z,t}=\mathcal{H}_{t}(i_{z,t}-\bar{\imath}_{z})+\bar{\imath}_{z}',
    '\mathcal{H}_{t}=\left( 1+\frac{i_{z,t}}{z_{t-1}}-\delta -g\right)^{\phi }',
    '\mathfrak{Gap}_{t}=\frac{1+r_{d}}{1+g}d_{t-1}-d_{t}+\frac{r_{dc,t-1}-g}{1+g}dc_{t-1}+\frac{r_{t-1}-g}{1+g}P_{t}b_{t-1} + P_{z,t}\mathbb{I}_{t}+\mathcal{T}_{o}-h_{o}P_{t}c_{t}-\mathcal{G}_{t}-\mathcal{N}_{t}-\mu {z^{e}}_{t-1}',
    '\mathfrak{Gap}_{t}=P_{t}\Delta b_{t}+\Delta d_{c,t}+(h_{t}-h_{o})P_{t}c_{t}-(\mathcal{T}_{t}-\mathcal{T}_{o})',
    'h_{t}^{\text{target}} = h_{o}+(1-\lambda )\frac{\mathfrak{Gap}_{t}}{P_{t}c_{t}}',
    '\mathcal{T}_{t}^{\text{target}}=\mathcal{T}_{o}-\lambda \mathfrak{Gap}_{t}',
    'h_{t}=Min\left\{ h_{t}^{r},h^{u}\right\}',
    '\mathcal{T}_{t}=Max\left\{ \mathcal{T}_{t}^{r},\mathcal{T}^{l}\right\}',
    'h_{t}^{r}=h_{t-1}+\lambda _{1}({\small h}_{t}^{\text{target}}-h_{t-1})+\lambda _{2}\frac{(x_{t-1}-x^{\text{target}})}{y_{t}} \quad \text{with} \quad \lambda _{1},\lambda _{2}>0',
    '\mathcal{T}_{t}^{r}=\mathcal{T}_{t-1}+\lambda _{3}(\mathcal{T}_{t}^{\text{target}}-\mathcal{T}_{t-1})-\lambda _{4}(x_{t-1}-x^{\text{target}}) \quad \text{with}\quad \lambda _{3},\lambda _{4}>0',
    'L_{x}+L_{n}=L',
    'q_{n,t}=\rho _{n}\left( \frac{P_{n,t}}{P_{t}}\right) ^{-\epsilon}c_{t}+a_{k}\left( i_{x,t}+i_{n,t}+AC_{x,t}+AC_{n,t}\right) +a_{z}\mathbb{I}_{z,t}',
    'd_{t}-d_{t-1}+d_{c,t}-d_{c,t-1}+{b^{\ast}}_{t}-{b^{\ast}}_{t-1} = \frac{r_{d}-g}{1+g}d_{t-1}+\frac{r_{dc,t-1}-g}{1+g}d_{c,t-1}+\frac{{r^{\ast}}_{t-1}-g}{1+g}{b^{\ast}}_{t-1} + \mathcal{P}_{t}+P_{z,t}\mathbb{I}_{z,t}+P_{k,t}\left(i_{x,t}+i_{n,t}+AC_{x,t}+AC_{n,t}\right) + P_{t}c_{t}-P_{n,t}q_{n,t}-P_{x,t}q_{x,t}-\mathcal{R}_{t}-\mathcal{G}_{t}-\mathcal{N}_{t}'
)

for ($i = 0; $i -lt $equationTexts.Length; $i++) {
    $row = $i + 2
    $wsEquations.Range("C$row").Value = $equationTexts[$i]
}

# --- Update sheet selections / active tab ---

# "params_deep" is no longer the active sheet; its selection moves to E6.
$wsParamsDeep.Activate()
$wsParamsDeep.Range("E6").Select()

# "equations" becomes the active sheet, with the selection on the last row (B45).
$wsEquations.Activate()
$wsEquations.Range("B45").Select()
